# Refresh the cached "datetimeFigureOut" date placeholder text across the
# slide master and every slide layout (Insert > Header & Footer > Date and
# time > Update automatically, refreshed on save) from 03.06.2025 to
# 11.11.2025.

$p = $ppt.ActivePresentation
$newDate = "11.11.2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
